$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.895781993865967
$ws.Range("B1").Value = 5.467911243438721
$ws.Range("C1").Value = 8.705711364746094
$ws.Range("D1").Value = 7.971055030822754
$ws.Range("E1").Value = 3.043112754821777
